# Auto-generated edit script: refresh market-price-derived Leve profit figures
# across all 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1089.9
$ws.Range("J17").Value = 1089.9
$ws.Range("L17").Value = 3269.7
$ws.Range("N17").Value = -3605.7
$ws.Range("H31").Value = 4150
$ws.Range("I31").Value = 4150
$ws.Range("K31").Value = 12450
$ws.Range("M31").Value = -12220
$ws.Range("H38").Value = 5845.5264
$ws.Range("I38").Value = 5845.5264
$ws.Range("K38").Value = 17536.5792
$ws.Range("M38").Value = -17164.5792
$ws.Range("H48").Value = 9009.5
$ws.Range("J48").Value = 9009.5
$ws.Range("L48").Value = 27028.5
$ws.Range("N48").Value = -27612.5
$ws.Range("H56").Value = 9009.5
$ws.Range("J56").Value = 9009.5
$ws.Range("L56").Value = 27028.5
$ws.Range("N56").Value = -28096.5
$ws.Range("H98").Value = 567.1667
$ws.Range("I98").Value = 567.1667
$ws.Range("K98").Value = 567.1667
$ws.Range("M98").Value = 930.8333
$ws.Range("H122").Value = 567.1667
$ws.Range("I122").Value = 567.1667
$ws.Range("K122").Value = 1701.5001
$ws.Range("M122").Value = 748.4999
$ws.Range("H127").Value = 3438.8
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H131").Value = 1512.6
$ws.Range("I131").Value = 1512.6
$ws.Range("K131").Value = 4537.799999999999
$ws.Range("M131").Value = 502.2000000000007
$ws.Range("H138").Value = 2960
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2960
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 8880
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -19160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 16110.333
$ws.Range("I122").Value = 16110.333
$ws.Range("K122").Value = 48330.999
$ws.Range("M122").Value = -45880.999
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470
$ws.Range("H140").Value = 60429
$ws.Range("J140").Value = 60429
$ws.Range("L140").Value = 60429
$ws.Range("N140").Value = -70789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1450
$ws.Range("I86").Value = 1450
$ws.Range("K86").Value = 1450
$ws.Range("M86").Value = -327
$ws.Range("H89").Value = 1450
$ws.Range("I89").Value = 1450
$ws.Range("K89").Value = 7250
$ws.Range("M89").Value = -1634
$ws.Range("H99").Value = 1833.3334
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H134").Value = 1083.5
$ws.Range("I134").Value = 704
$ws.Range("K134").Value = 2112
$ws.Range("M134").Value = 423

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2474.3333
$ws.Range("I99").Value = 2474.3333
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2474.3333
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -976.3332999999998
$ws.Range("N99").ClearContents()
$ws.Range("H122").Value = 3499.5
$ws.Range("I122").Value = 3499.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10498.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8048.5
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2474.3333
$ws.Range("I126").Value = 2474.3333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7422.999899999999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4952.999899999999
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 7907.3335
$ws.Range("I132").Value = 4497
$ws.Range("K132").Value = 13491
$ws.Range("M132").Value = -10961
$ws.Range("H134").Value = 4609
$ws.Range("I134").Value = 1350
$ws.Range("J134").Value = 5333.222
$ws.Range("K134").Value = 4050
$ws.Range("L134").Value = 15999.666
$ws.Range("M134").Value = -1515
$ws.Range("N134").Value = -21069.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 244.75
$ws.Range("J46").Value = 486.5
$ws.Range("L46").Value = 1459.5
$ws.Range("N46").Value = -1641.5
$ws.Range("H49").Value = 5000
$ws.Range("J49").Value = 5000
$ws.Range("L49").Value = 15000
$ws.Range("N49").Value = -15312
$ws.Range("H92").Value = 1995
$ws.Range("J92").Value = 1995
$ws.Range("L92").Value = 5985
$ws.Range("N92").Value = -8481
$ws.Range("H121").Value = 1537.6
$ws.Range("I121").Value = 1537.6
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 4612.799999999999
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -3302.799999999999
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1501.6666
$ws.Range("I102").Value = 1170
$ws.Range("K102").Value = 1170
$ws.Range("M102").Value = 452
$ws.Range("H132").Value = 2204
$ws.Range("I132").Value = 2204
$ws.Range("K132").Value = 6612
$ws.Range("M132").Value = -4082

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 416.66666
$ws.Range("H27").Value = 416.66666
$ws.Range("H40").Value = 4202
$ws.Range("I40").Value = 4202
$ws.Range("K40").Value = 4202
$ws.Range("M40").Value = -4066
$ws.Range("H46").Value = 2002
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2002
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2002
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2378
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H136").Value = 14749.75
$ws.Range("I136").Value = 14749.75
$ws.Range("K136").Value = 44249.25
$ws.Range("M136").Value = -41699.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 49390
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H41").Value = 29999.5
$ws.Range("I41").Value = 29999.5
$ws.Range("K41").Value = 29999.5
$ws.Range("M41").Value = -29609.5
$ws.Range("H56").Value = 20314
$ws.Range("J56").Value = 20314
$ws.Range("L56").Value = 20314
$ws.Range("N56").Value = -21742
$ws.Range("H88").Value = 5171
$ws.Range("I88").Value = 5171
$ws.Range("K88").Value = 5171
$ws.Range("M88").Value = -4765
$ws.Range("H91").Value = 5171
$ws.Range("I91").Value = 5171
$ws.Range("K91").Value = 5171
$ws.Range("M91").Value = -3767
$ws.Range("H123").Value = 32250
$ws.Range("J123").Value = 32250
$ws.Range("L123").Value = 32250
$ws.Range("N123").Value = -42050
$ws.Range("H132").Value = 4422
$ws.Range("I132").Value = 4422
$ws.Range("K132").Value = 13266
$ws.Range("M132").Value = -10736
